# Ingredient_recette.xlsx edit script
# - Fix "kg" -> "Kg" for the unite column on the MPSLPO0001 row
# - Add a new recipe block: REC_VACHE_BRASSE_SUCRE with two ingredients
#   (BASE_VACHE_BRASSE_NATURE and MPSING0001)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the leftover no-op formatting that used to sit on every cell / column
# (applyAlignment="1" with no actual alignment) - restores the default style.
$ws.Cells.ClearFormats()

# Re-apply the real formatting that the cells in column B carry
$ws.Range("B2:B3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = -4108

# Correct the casing of the unit on row 3 (ingredient MPSLPO0001)
$ws.Range("D3").Value = "Kg"

# New rows for the REC_VACHE_BRASSE_SUCRE recipe
# (row 6 entered before row 5 to match the original shared-string order)
$ws.Range("A6").Value = "REC_VACHE_BRASSE_SUCRE"
$ws.Range("B6").Value = "MPSING0001"
$ws.Range("C6").Value = 60
$ws.Range("D6").Value = "Kg"
$ws.Range("E6").Value = "True"

$ws.Range("A5").Value = "REC_VACHE_BRASSE_SUCRE"
$ws.Range("B5").Value = "BASE_VACHE_BRASSE_NATURE"
$ws.Range("C5").Value = 940
$ws.Range("D5").Value = "L"
$ws.Range("E5").Value = "True"

# Selection moved to C7 in the saved file
$ws.Range("C7").Select()
